# Actualización 11 de Mayo - Mañana
# Rescatables sheet: add a new rescatable record for HERNANDEZ DOLORES JOEL
# EDUARDO ahead of the existing FLORES DE LOS SANTOS JHOVANA record, and
# correct the latter's "Reprobadas" count from 2 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Push the existing row 2 down to row 3, then drop the formatting that
# Insert() copies down from the header row so the new row matches the
# plain (unstyled) data rows used elsewhere in the sheet.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New student record (row 2)
$ws.Range("A2").Value = 18330051920248
$ws.Range("B2").Value = "HERNANDEZ"
$ws.Range("C2").Value = "DOLORES"
$ws.Range("D2").Value = "JOEL EDUARDO"
$ws.Range("E2").Value = "ANALIZA SANGRE MEDIANTE PRUEBAS HORMONALES, TOXICOLÓGICAS Y DE MARCADORES TUMORALES"
$ws.Range("F2").Value = "6ALCV"
$ws.Range("G2").Value = 2

# Existing student record, now on row 3 - only the Reprobadas count changes.
$ws.Range("G3").Value = 1
